$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("C4").Value = 2.844
$ws.Range("D4").Value = 3.105
$ws.Range("E4").Value = 2.388
$ws.Range("F4").Value = 2.464
$ws.Range("G4").Value = 2.141
$ws.Range("H4").Value = 3.269

# Row 5
$ws.Range("C5").Value = 2.647
$ws.Range("D5").Value = 3.098
$ws.Range("E5").Value = 2.37
$ws.Range("F5").Value = 2.809
$ws.Range("G5").Value = 2.012
$ws.Range("H5").Value = 3.252

# Row 6
$ws.Range("C6").Value = 0.752
$ws.Range("D6").Value = 0.461
$ws.Range("E6").Value = 0.443
$ws.Range("F6").Value = 0.648
$ws.Range("G6").Value = 0.488
$ws.Range("H6").Value = 0.428

# Row 7
$ws.Range("C7").Value = 0.966
$ws.Range("D7").Value = 0.458
$ws.Range("E7").Value = 0.438
$ws.Range("F7").Value = 0.728
$ws.Range("G7").Value = 0.485
$ws.Range("H7").Value = 0.426
